# Adds "Description" and "Acceptance Criteria" support to the ticket
# creation sheet, replaces the sample data with the coffee-machine /
# meeting-room dummy tickets, and tidies up the "info" helper sheet.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("info")

# ---------------------------------------------------------------------
# Sheet1 ("tickets to create")
# ---------------------------------------------------------------------

# New header row: Issue Type | Project | Summary | Description | Acceptance Criteria
$ws1.Cells.Item(1,1).Value = "Issue Type"
$ws1.Cells.Item(1,2).Value = "Project"
$ws1.Cells.Item(1,3).Value = "Summary"
$ws1.Cells.Item(1,4).Value = "Description"
$ws1.Cells.Item(1,5).Value = "Acceptance Criteria - User Stories Only"

# Bold the header row (A1:C1 plain bold, D1:E1 bold + wrap)
$ws1.Range("A1:E1").Font.Bold = $true
$ws1.Range("D1:E1").WrapText = $true

# Body rows: Epic + three user stories
$ws1.Cells.Item(2,1).Value = 5
$ws1.Cells.Item(2,2).Value = "GRW"
$ws1.Cells.Item(2,3).Value = "DUMMY EPIC: Fix Coffee Machine"
$ws1.Cells.Item(2,4).Value = "As a developer, I've noticed the office coffee machine now requires a secret handshake before dispensing caffeine. "
$ws1.Cells.Item(2,5).Value = "Machine should dispense coffee without existential questioning. We need either an IT exorcism or a better espresso setup."

$ws1.Cells.Item(3,1).Value = 6
$ws1.Cells.Item(3,2).Value = "GRW"
$ws1.Cells.Item(3,3).Value = "DUMMY TICKET: Beverage Temperature Retention Enhancement"
$ws1.Cells.Item(3,4).Value = "As a consultant, I want my coffee to stay hot for longer than 10 minutes"
$ws1.Cells.Item(3,5).Value = "Coffee remains at 140°F+ for at least 30 minutes after brewing"

$ws1.Cells.Item(4,1).Value = 6
$ws1.Cells.Item(4,2).Value = "GRW"
$ws1.Cells.Item(4,3).Value = "DUMMY TICKET: Meeting Room Seating Capacity Audit & Restoration"
$ws1.Cells.Item(4,4).Value = "As a consultant, I want meeting rooms to actually have chairs"
$ws1.Cells.Item(4,5).Value = "All meeting rooms have minimum 8 functional chairs rated for 8-hour workdays"

$ws1.Cells.Item(5,1).Value = 6
$ws1.Cells.Item(5,2).Value = "GRW"
$ws1.Cells.Item(5,3).Value = "DUMMY TICKET: Establish Conference Call Audio Clarity Protocol"
$ws1.Cells.Item(5,4).Value = "As a consultant, I want to stop asking 'Can everyone hear me?' five times per call"
$ws1.Cells.Item(5,5).Value = 'Echo cancellation reduces repetitive "Can you hear me now?" questions by 90%'

# Wrap the new Description / Acceptance Criteria columns for all data rows
$ws1.Range("D2:E5").WrapText = $true

# Give each data row enough height to show the wrapped text
$ws1.Rows.Item(2).RowHeight = 57.6
$ws1.Rows.Item(3).RowHeight = 43.2
$ws1.Rows.Item(4).RowHeight = 43.2
$ws1.Rows.Item(5).RowHeight = 43.2

# Column widths
$ws1.Columns.Item(1).ColumnWidth = 8.944010416666666
$ws1.Columns.Item(2).ColumnWidth = 6.166666666666667
$ws1.Columns.Item(3).ColumnWidth = 57.721354166666664
$ws1.Columns.Item(4).ColumnWidth = 28.053385416666668
$ws1.Columns.Item(5).ColumnWidth = 32.830729166666664

# Print orientation (portrait)
$ws1.PageSetup.Orientation = 1

# Selection shown when the sheet is active
[void]$ws1.Range("C3").Select()

# ---------------------------------------------------------------------
# "info" helper sheet - the "<- Most common" note moves from the Task
# row down to the User Story row.
# ---------------------------------------------------------------------
$ws2.Cells.Item(3,3).ClearContents()
$ws2.Cells.Item(4,3).Value = "← Most common"

[void]$ws2.Range("G4").Select()

# Sheet1 must stay the active sheet/tab when the workbook is reopened.
[void]$ws1.Activate()
[void]$ws1.Range("C3").Select()
